$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H4").Value = "no 4:40pm-5:40pm"
$ws.Range("H8").Value = "no 4:40pm-5:40pm"
$ws.Range("H5").Value = "Only After 12:00pm"
$ws.Range("H6").Value = "no 1/12/14-1/20/14, no 3/3/14"
$ws.Range("H7").Value = "no 2/17/13"
$ws.Range("H9").Value = "only thursday, no saturday"

$ws.Range("H6").Select()

$ws.Rows.Item(4).RowHeight = 13.45
$ws.Rows.Item(5).RowHeight = 13.45
$ws.Rows.Item(6).RowHeight = 13.45
$ws.Rows.Item(7).RowHeight = 13.45
$ws.Rows.Item(9).RowHeight = 14.45
